$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.196.61'
$ws.Range('E2').Value = '  +11.22%  '
$ws.Range('D3').Value = '1.820.12'
$ws.Range('E3').Value = '  +7.55%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.544'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.09%  '
$ws.Range('E7').Value = '  +0.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.24'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.43%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.97'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.281'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0673'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0931'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.90%  '
$ws.Range('D13').Value = '2.089.20'
$ws.Range('E13').Value = '  +7.82%  '
$ws.Range('D14').Value = '1.822.92'
$ws.Range('E14').Value = '  +7.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.643'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.13%  '
$ws.Range('D16').Value = '34.190.67'
$ws.Range('E16').Value = '  +11.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '10.25'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.53%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.33'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +7.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.74'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '258.08'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.58%  '
$ws.Range('D21').Value = '0.0₃0750'
$ws.Range('E21').Value = '  +4.07%  '
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.47'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.97%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.37'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.76%  '
$ws.Range('E25').Value = '  +0.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.60'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.15'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.94%  '
$ws.Range('E29').Value = '  +4.00%  '
$ws.Range('E30').Value = '  +0.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.87'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +11.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0516'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.95%  '
$ws.Range('E33').Value = '  +6.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.56'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.56%  '
$ws.Range('D35').Value = '1.573.25'
$ws.Range('E35').Value = '  +3.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.84'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.07'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.53%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0189'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.10%  '
$ws.Range('B39').Value = 'Aave'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '84.87'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.79%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.627'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.85'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.43%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.34'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.913'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.13'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0521'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.12%  '
$ws.Range('E46').Value = '  +4.50%  '
$ws.Range('E47').Value = '  +8.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.74'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.37%  '
$ws.Range('E49').Value = '  +0.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.84'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.75%  '
$ws.Range('E51').Value = '  +5.95%  '
